# SA_Info.xlsx - add 7 new XPath entries to the "XPath" sheet (rows 48-54)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XPath")
$ws.Activate()

# Green, bold, 12pt Menlo used for every xpath value in column B (matches
# the formatting already used for rows 45-47 just above the new block).
$xpathFontName  = "Menlo"
$xpathFontSize  = 12
$xpathFontColor = 225295   # RGB(15,112,3) == FF0F7003 packed as BGR int

function Set-LabelCell($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.WrapText = $true
    $c.Font.Name = "Arial"
    $c.Font.Size = 10
    $c.Font.Bold = $false
}

function Set-XPathCell($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.WrapText = $true
    $c.Font.Bold = $true
    $c.Font.Size = $xpathFontSize
    $c.Font.Color = $xpathFontColor
    $c.Font.Name = $xpathFontName
}

$rows = @(
    @{ Row = 48; Label = "ordered on";          XPath = "(//SPAN[text()='28 Jan 18, 10:43 PM'][text()='28 Jan 18, 10:43 PM'])[2]"; Height = 15 },
    @{ Row = 49; Label = "total price";          XPath = "(//SPAN[text()='14599'])[3]"; Height = 15 },
    @{ Row = 50; Label = "channel";              XPath = "(//SPAN[text()='AndroidApp'])[1]"; Height = 15 },
    @{ Row = 51; Label = "order verification";   XPath = "(//DIV[@class='_-sa-flipkart-src-Components-DT-DTButton-DTButton_button _-sa-flipkart-src-Components-DT-DTButton-DTButton_outlineBtn _-sa-flipkart-src-Components-DT-DTButton-DTButton_isDisabled'])[1]"; Height = 28.2 },
    @{ Row = 52; Label = "create incident";      XPath = "(//DIV[@class='_-sa-flipkart-src-Components-DT-DTButton-DTButton_button _-sa-flipkart-src-Components-DT-DTButton-DTButton_outlineBtn _-sa-flipkart-src-Components-DT-DTButton-DTButton_isDisabled'])[2]"; Height = 28.2 },
    @{ Row = 53; Label = "price adjustment";     XPath = "(//DIV[@class='_-sa-flipkart-src-Components-DT-DTButton-DTButton_button _-sa-flipkart-src-Components-DT-DTButton-DTButton_outlineBtn _-sa-flipkart-src-Components-DT-DTButton-DTButton_isDisabled'])[3]"; Height = 28.2 },
    @{ Row = 54; Label = "cancel dt";            XPath = "(//DIV[@class='_-sa-flipkart-src-Components-DT-DTButton-DTButton_button _-sa-flipkart-src-Components-DT-DTButton-DTButton_outlineBtn _-sa-flipkart-src-Components-DT-DTButton-DTButton_isDisabled'])[4]"; Height = 28.2 }
)

foreach ($r in $rows) {
    Set-LabelCell ("A" + $r.Row) $r.Label
    Set-XPathCell ("B" + $r.Row) $r.XPath
    $ws.Rows.Item($r.Row).RowHeight = $r.Height
}

# Move the view roughly the same amount it moved in the source edit (the
# sheet scrolled down a few rows and the active cell followed the newly
# typed data down to B59).
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B59").Select()
